$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "autoawq" / "A" row), shifting all rows below up by one.
$ws.Rows.Item(2).Delete()

# Update the selection to match the post-edit state (row 2, full row selection).
$ws.Range("A2:XFD2").Select()
